# switcherr method integrated with gpsSearcher
#
# Adds the new European cities gathered by the gpsSearcher/"switcherr"
# pass to Sheet2 (the running log) and refreshes Sheet1 (the "current /
# latest" snapshot sheet) so that it shows the most recently looked-up
# trip (Madrid, bilbao, valencia) instead of the older
# st petersberg / helsinki entries - those two rows move down into the
# Sheet2 log instead.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$monday = "Monday the 31st August"

# ---------------------------------------------------------------------
# Sheet2: append the rows that were produced by this search pass.
# Rows 15-16 restate the previous st petersberg / helsinki find (they
# used to live on Sheet1), rows 17-31 are the brand new cities.
# ---------------------------------------------------------------------
$sheet2Rows = @(
    @("st petersberg", "59.9311° N, 30.3609° E", $monday),
    @("helsinki",      "60.1699° N, 24.9384° E", $monday),
    @("oslo",          "59.9139° N, 10.7522° E", $monday),
    @("malmo",         "55.6050° N, 13.0038° E", $monday),
    @("copenhagen",    "55.6761° N, 12.5683° E", $monday),
    @("odense",        "55.4038° N, 10.4024° E", $monday),
    @("kiel",          "54.3233° N, 10.1228° E", $monday),
    @("rostock",       "54.0924° N, 12.0991° E", $monday),
    @("Svalbard",      "77.8750° N, 20.9752° E", $monday),
    @("Tromso",        "69.6492° N, 18.9553° E", $monday),
    @("Hannover",      "52.3759° N, 9.7320° E",  $monday),
    @("Hamburg",       "53.5511° N, 9.9937° E",  $monday),
    @("London",        "51.5074° N, 0.1278° W",  $monday),
    @("Lisbon",        "38.7223° N, 9.1393° W",  $monday),
    @("Madrid",        "40.4168° N, 3.7038° W",  $monday),
    @("bilbao",        "43.2630° N, 2.9350° W",  $monday),
    @("valencia",      "39.4699° N, 0.3763° W",  $monday)
)

$startRow = 15
for ($i = 0; $i -lt $sheet2Rows.Count; $i++) {
    $row = $startRow + $i
    $vals = $sheet2Rows[$i]
    $ws2.Cells.Item($row, 1).Value = $vals[0]
    $ws2.Cells.Item($row, 2).Value = $vals[1]
    $ws2.Cells.Item($row, 3).Value = $vals[2]
}

# ---------------------------------------------------------------------
# Sheet1: swap in the newest trip (Madrid / bilbao / valencia), which
# now occupies 3 rows instead of the previous 2.
# ---------------------------------------------------------------------
$sheet1Rows = @(
    @("Madrid",   "40.4168° N, 3.7038° W"),
    @("bilbao",   "43.2630° N, 2.9350° W"),
    @("valencia", "39.4699° N, 0.3763° W")
)

for ($i = 0; $i -lt $sheet1Rows.Count; $i++) {
    $row = 2 + $i
    $vals = $sheet1Rows[$i]
    $ws1.Cells.Item($row, 1).Value = $vals[0]
    $ws1.Cells.Item($row, 2).Value = $vals[1]
}

# ---------------------------------------------------------------------
# Selections: leave a transient selection on Sheet1 (A4, the new last
# row) but finish on Sheet2 (C22:C31) so Sheet2 stays the active tab.
# ---------------------------------------------------------------------
$ws1.Range("A4").Select() | Out-Null
$ws2.Select() | Out-Null
$ws2.Range("C22:C31").Select() | Out-Null
